$d = $word.ActiveDocument

# Replace responsible-technician name (appears in the title block and in the signature block)
$d.Content.Find.Execute("RENAN NUNES ZERINO", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ELENILDA FERREIRA", 2)

# Replace course completion date
$d.Content.Find.Execute("30/04/2025", $true, $false, $false, $false, $false,
                         $true, 1, $false, "10/07/2025", 2)

# Replace course/student ID code
$d.Content.Find.Execute("123456:654321", $true, $false, $false, $false, $false,
                         $true, 1, $false, "316595:25101", 2)
